# Saldo.xlsx update — refresh the "Export" balances snapshot.
#
# Net effect (top of the sheet, below the header + first CARLOS row):
#   - Rows for AYRTON / CHRISTIAN / CESAR / GUSTAVO(7573.39) / NILBORN(7573.39)
#     are replaced by three new rows: ANDREA, TATIANA and MAGALI (with a new
#     13693.23 balance).
#   - The CARLOS(5306.81) / RAPHAELA(5200) rows are removed entirely.
#   - The old MAGALI(95.47) row further down is removed entirely (MAGALI now
#     only appears once, near the top, with the new balance).
#
# We delete rows bottom-up so earlier row numbers stay valid while we work,
# then overwrite the three rows that take on new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old, lower "MAGALI" row (originally row 34).
$ws.Rows.Item(34).Delete()

# 2) Remove the CARLOS(5306.81) and RAPHAELA(5200) rows (originally rows 9-10).
$ws.Rows("9:10").Delete()

# 3) Remove two of the five rows being collapsed down to three
#    (GUSTAVO 7573.39 and NILBORN 7573.39 — originally rows 6-7),
#    leaving rows 3-5 as the three rows to overwrite.
$ws.Rows("6:7").Delete()

# 4) Overwrite rows 3-5 with the new account/name/balance data.
#    A leading apostrophe forces text entry (keeping the leading zeros in
#    the account numbers); ClearFormats() then drops the implicit "Text"
#    number-format Excel would otherwise stamp on the cell, so the cell
#    keeps the plain/default style it had before (matches the other
#    untouched account-number cells in the sheet).
$ws.Range("A3").Value = "'004181486"
$ws.Range("A3").ClearFormats()
$ws.Range("B3").Value = "ANDREA"
$ws.Range("C3").Value = 32392.06

$ws.Range("A4").Value = "'005366671"
$ws.Range("A4").ClearFormats()
$ws.Range("B4").Value = "TATIANA"
$ws.Range("C4").Value = 13708

$ws.Range("A5").Value = "'004207641"
$ws.Range("A5").ClearFormats()
$ws.Range("B5").Value = "MAGALI"
$ws.Range("C5").Value = 13693.23
